$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire "Group job migration" column (column E), shifting
# the "Statut" column (old F) left into its place.
$ws.Columns.Item(5).EntireColumn.Delete()

# Update the selected cell on the sheet view to C21
$ws.Range("C21").Select()
